# Apply the commit's content changes to the workbook.
$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: update URL, Version, Date and Publisher values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/reinsurance-met-indicator"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- "Elements" sheet: the ele-1/ext-1 constraint note was moved off the root
#     "Extension" row (row 2) - it already correctly appears on the
#     "Extension.extension" row (row 4), so only row 2 needs clearing. ---
$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").Value = ""

$wb.Save()
